$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    # Prefix with an apostrophe so numeric-looking strings stay text,
    # then reset style so no stray "quote prefix" number format sticks
    # around on the cell (keeps formatting identical to the original).
    $cell.Value = "'" + $newValue
    $cell.Style = "Normal"
}

# Price column (D) updates
Set-TextValue "D2" "247.62"
Set-TextValue "D4" "5.292"
Set-TextValue "D5" "0.05723"
Set-TextValue "D6" "3.440"
Set-TextValue "D7" "0.8100"
Set-TextValue "D8" "0.8755"
Set-TextValue "D9" "0.1429"
Set-TextValue "D10" "0.07384"
Set-TextValue "D12" "0.03128"
Set-TextValue "D13" "0.09402"
Set-TextValue "D14" "3.985"
Set-TextValue "D15" "0.001575"
Set-TextValue "D16" "0.04830"
Set-TextValue "D17" "0.0005852"
Set-TextValue "D18" "0.006145"
Set-TextValue "D19" "0.005123"
Set-TextValue "D20" "0.0009958"
Set-TextValue "D22" "3.735"
Set-TextValue "D23" "6.320"
Set-TextValue "D25" "0.3279"
Set-TextValue "D40" "0.03934"
Set-TextValue "D41" "0.006761"
Set-TextValue "D42" "0.1068"
Set-TextValue "D43" "0.002610"
Set-TextValue "D44" "0.007505"
Set-TextValue "D45" "0.00005620"
Set-TextValue "D47" "0.6002"
Set-TextValue "D49" "0.00002101"

# Volume(1h) column (E) updates
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
